$d = $word.ActiveDocument

$r2 = $d.Range(96,117)
$found = $r2.Find.Execute("Proposa_V1.0", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Proposa_V1.2", 2)
Write-Output "found=$found r2 start=$($r2.Start) end=$($r2.End) text=[$($r2.Text)]"

$bm = $d.Range(117, 117)
Write-Output "bm start=$($bm.Start) end=$($bm.End) text=[$($bm.Text)]"
$d.Bookmarks.Add("_GoBack", $bm)
Write-Output "added"
Write-Output "bookmarks count = $($d.Bookmarks.Count)"
$b = $d.Bookmarks(1)
Write-Output "bookmark1 name=$($b.Name) start=$($b.Range.Start) end=$($b.Range.End)"
